$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "To check that it is possible to read the temperature ..." paragraph:
#    insert a new "sudo cat " run right before the
#    "/sys/bus/w1/devices/28-f3a49d1964ff/w1_slave" text.
# ------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "at the command line enter '/sys/bus/w1/devices/28-f3a49d1964ff/w1_slave",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "at the command line enter 'sudo cat /sys/bus/w1/devices/28-f3a49d1964ff/w1_slave",
    2)
Write-Host "Step1 (sudo cat insert) found/replaced:" $found1

# ------------------------------------------------------------------
# 2) "To check that the new nodes have been detected, select ..." paragraph:
#    collapse the runs that spell out
#    ", select 'Node and Zone Configuration' from the Settings dropdown
#    list and click the 'Nodes' button." into a single run (text itself
#    is unchanged).
# ------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "To check that the new nodes have been detected, select 'Node and Zone Configuration' from the Settings dropdown list and click the 'Nodes' button.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To check that the new nodes have been detected, select 'Node and Zone Configuration' from the Settings dropdown list and click the 'Nodes' button.",
    2)
Write-Host "Step2 (Node and Zone merge) found/replaced:" $found2

# ------------------------------------------------------------------
# 3) "Select 'Device Configuration' ..." paragraph:
#    collapse the runs that spell out
#    "Select 'Device Configuration' from the Settings dropdown list and
#    click the 'Sensors' button." into a single run (text itself is
#    unchanged).
# ------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(
    "Select 'Device Configuration' from the Settings dropdown list and click the 'Sensors' button.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Select 'Device Configuration' from the Settings dropdown list and click the 'Sensors' button.",
    2)
Write-Host "Step3 (Device Configuration merge) found/replaced:" $found3
